# Carnet de bord - Stage semaine 3 - edits
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1) Insert one new row inside the "VENDREDI" block (between the old rows
#    29 and 30) so the block grows from 4 rows (28-31) to 5 rows (28-32).
#    Everything below (Commentaires, legend, ...) shifts down by one row,
#    which is exactly what the target layout needs.
# ---------------------------------------------------------------------------
$ws.Rows.Item(30).Insert()

# The "VENDREDI après-midi" header cell used to be the merged range A30:A31;
# after the insert it is now A31:A32. Re-merge it as A30:A32 so the header
# spans the whole afternoon block (3 rows) like in the target file.
$apresMidiLabel = $ws.Cells.Item(31,1).Value2
$ws.Range("A31:A32").UnMerge()
$ws.Cells.Item(31,1).Value2 = ""
$ws.Cells.Item(30,1).Value2 = $apresMidiLabel
$ws.Range("A30:A32").Merge()
$ws.Range("A30:A32").HorizontalAlignment = -4108
$ws.Range("A30:A32").VerticalAlignment = -4108
$ws.Range("A30:A32").WrapText = $true
$ws.Rows.Item(30).RowHeight = 18

# ---------------------------------------------------------------------------
# 2) Fill in the "VENDREDI matin" missions (rows 28-29/30)
# ---------------------------------------------------------------------------

# Mission 1 (single row 28): "Transfert de données disquettes vers CD"
$ws.Cells.Item(28,2).Value2 = "Transfert de données disquettes vers CD"
$ws.Range("C28:D28").Value2 = "Transfert de données disquettes vers CD"
$ws.Cells.Item(28,5).Value2 = "Logiciel Nero Express"
$ws.Cells.Item(28,6).Value2 = "30 min"
$ws.Cells.Item(28,7).Value2 = "Oberserver et manipuler avec monsieur SEGATO"
$ws.Cells.Item(28,8).Value2 = 2
$ws.Cells.Item(28,9).Value2 = "Etre autonome ef faire des corrélations avec mes connaissances du stage"
$ws.Cells.Item(28,10).Value2 = 4
$ws.Rows.Item(28).RowHeight = 39.6
$ws.Range("A28:J28").WrapText = $true

# Mission 2 (merged rows 29-30): "Gérer la clientèle" - spans matin/après-midi
$ws.Range("B29:B30").Merge()
$ws.Range("C29:D30").Merge()
$ws.Range("E29:E30").Merge()
$ws.Range("F29:F30").Merge()
$ws.Range("G29:G30").Merge()
$ws.Range("H29:H30").Merge()
$ws.Range("I29:I30").Merge()
$ws.Range("J29:J30").Merge()

$ws.Cells.Item(29,2).Value2 = "Gérer la clientèle"
$ws.Cells.Item(29,3).Value2 = "Gérer les demandes de la clientèle ainsi que certaines factures et devis"
$ws.Cells.Item(29,5).Value2 = "Poste informatique de l'accueil, imprimante"
$ws.Cells.Item(29,6).Value2 = "Toute la journée"
$ws.Cells.Item(29,7).Value2 = "Savoir communiquer avec le client et être autonome (ou poser des questions à M. SEGATO)"
$ws.Cells.Item(29,8).Value2 = "2 ou 3"
$ws.Cells.Item(29,9).Value2 = "Oberserver, questionner, faire des corrélations avec mes connaissances du stage"
$ws.Cells.Item(29,10).Value2 = 3
$ws.Range("B29:J30").WrapText = $true
$ws.Range("B29:J30").HorizontalAlignment = -4108
$ws.Range("B29:J30").VerticalAlignment = -4108

# ---------------------------------------------------------------------------
# 3) Fill in the "VENDREDI après-midi" missions (rows 31 and 32)
# ---------------------------------------------------------------------------

# Mission 3 (row 31): "Retrait marchandise"
$ws.Cells.Item(31,2).Value2 = "Retrait marchandise"
$ws.Range("C31:D31").Value2 = "Retrait d'un PC fixe, de 3 PC portables et d'un adaptateur pour carte mémoire chez le fournisseur NEO"
$ws.Cells.Item(31,5).Value2 = "Voiture personnelle"
$ws.Cells.Item(31,6).Value2 = "1h"
$ws.Cells.Item(31,7).Value2 = "Savoir aller à une adresse demandé et récupérer la marchandise"
$ws.Cells.Item(31,8).Value2 = 3
$ws.Cells.Item(31,9).Value2 = "Savoir aller à une adresse demandé et récupérer la marchandise"
$ws.Cells.Item(31,10).Value2 = 3
$ws.Rows.Item(31).RowHeight = 45
$ws.Range("A31:J31").WrapText = $true

# Mission 4 (row 32): "Test PC portable"
$ws.Cells.Item(32,2).Value2 = "Test PC portable"
$ws.Range("C32:D32").Value2 = "Test de la batterie d'un PC portable pour voir si il fonctionne et le rendre au client après réparation "
$ws.Cells.Item(32,5).Value2 = "Atelier de l'entreprise avec le matériel informatique nécessaire"
$ws.Cells.Item(32,6).Value2 = "1h"
$ws.Cells.Item(32,7).Value2 = "Savoir faire des manipulations simples que l'on m'a montré"
$ws.Cells.Item(32,8).Value2 = 3
$ws.Cells.Item(32,9).Value2 = "Etre autonome ef faire des corrélations avec mes connaissances du stage"
$ws.Cells.Item(32,10).Value2 = 4
$ws.Rows.Item(32).RowHeight = 33
$ws.Range("A32:J32").WrapText = $true

# ---------------------------------------------------------------------------
# 4) Comments section (rows 33-36 after the shift)
# ---------------------------------------------------------------------------
# Row 34 (previously the "Anthony" row) gets the surname added.
$ws.Cells.Item(34,1).Value2 = "J'ai été accompagné d'un autre stagiaire lycéen cette semaine : Anthony CAVAGNÉ."

# Row 35 was an empty comment line; fill it in with the new end-of-week note.
$ws.Cells.Item(35,1).Value2 = "Cette semaine s'est bien déroulé malgré le fait que je n'ai pas eu l'occassion de demander  à monsieur SEGATO de remplir les documents de fin de stage ainsi que de le remercier pour le stage car il était en intervention toute la journée."
$ws.Range("A35:J35").WrapText = $true

# ---------------------------------------------------------------------------
# 5) Update the print area to the new sheet extent
# ---------------------------------------------------------------------------
$printAreaName = $wb.Names.Item("carnet de bord!Print_Area")
$printAreaName.RefersTo = "='carnet de bord'!`$A`$1:`$J`$42"

Write-Host "Done"
